$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D (media) and column E (desviacion) for rows 2-7
$ws.Range("D2").Value = 9.26666666666667
$ws.Range("E2").Value = 1.09937108386811

$ws.Range("D3").Value = 8.875
$ws.Range("E3").Value = 1.18074122128849

$ws.Range("D4").Value = 8.69
$ws.Range("E4").Value = 1.69388879778433

$ws.Range("D5").Value = 7.5
$ws.Range("E5").Value = 2.27751505591774

$ws.Range("D6").Value = 7.39166666666667
$ws.Range("E6").Value = 2.20210273452969

$ws.Range("D7").Value = 7.02833333333333
$ws.Range("E7").Value = 2.65052511579311
